# Applies the stock-report correction described in the commit diff.
# For each affected item row: Quantity (F) is adjusted and Value (G = D * F)
# is recomputed; for a handful of rows the Item Code (B), Selling Price (E),
# Quantity (F) and Value (G) are swapped between two adjacent duplicate-name
# rows; Sub Total / Grand Total rows (B only) are updated to match the new
# sum of the Value column within their section.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F71").Value = 368
$ws.Range("G71").Value = 23441.6

$ws.Range("F77").Value = 287
$ws.Range("G77").Value = 13414.38

$ws.Range("F86").Value = 84
$ws.Range("G86").Value = 10539.48

$ws.Range("B90").Value = 202074.62

$ws.Range("F115").Value = 231
$ws.Range("G115").Value = 22363.11

$ws.Range("B117").Value = 16415.39

$ws.Range("F144").Value = 1226
$ws.Range("G144").Value = 10359.7

$ws.Range("F145").Value = 658
$ws.Range("G145").Value = 5257.42

$ws.Range("B147").Value = 18647.96

$ws.Range("F149").Value = 253
$ws.Range("G149").Value = 16394.4

$ws.Range("F151").Value = 104
$ws.Range("G151").Value = 9035.52

$ws.Range("F152").Value = 73
$ws.Range("G152").Value = 6445.17

$ws.Range("B156").Value = 35990.32

$ws.Range("F163").Value = 14
$ws.Range("G163").Value = 3723.44

$ws.Range("F164").Value = 72
$ws.Range("G164").Value = 8253.360000000001

$ws.Range("B175").Value = 33407.72

$ws.Range("F197").Value = 24
$ws.Range("G197").Value = 1489.44

$ws.Range("F203").Value = 70
$ws.Range("G203").Value = 1411.2

$ws.Range("F205").Value = 31
$ws.Range("G205").Value = 11691.34

$ws.Range("F211").Value = 55
$ws.Range("G211").Value = 5566

$ws.Range("B216").Value = 48943.91

$ws.Range("B227").Value = 55373
$ws.Range("E227").Value = 163.62
$ws.Range("F227").Value = -94
$ws.Range("G227").Value = -13562.32

$ws.Range("B228").Value = 63520
$ws.Range("E228").Value = 153.4
$ws.Range("F228").Value = 67
$ws.Range("G228").Value = 9666.76

$ws.Range("B229").Value = 63531
$ws.Range("E229").Value = 152.53
$ws.Range("F229").Value = 67
$ws.Range("G229").Value = 9613.16

$ws.Range("B230").Value = 57802
$ws.Range("E230").Value = 162.71
$ws.Range("F230").Value = -79
$ws.Range("G230").Value = -11334.92

$ws.Range("F234").Value = 43
$ws.Range("G234").Value = 2206.76

$ws.Range("F256").Value = 294
$ws.Range("G256").Value = 44443.98

$ws.Range("B260").Value = 210595.72

$ws.Range("F270").Value = 42
$ws.Range("G270").Value = 1354.08

$ws.Range("B275").Value = 7668.28

$ws.Range("F280").Value = 147
$ws.Range("G280").Value = 24863.58

$ws.Range("F282").Value = 11
$ws.Range("G282").Value = 590.7

$ws.Range("F283").Value = 47
$ws.Range("G283").Value = 16049.09

$ws.Range("F285").Value = 16
$ws.Range("G285").Value = 446.88

$ws.Range("F294").Value = 49
$ws.Range("G294").Value = 3496.64

$ws.Range("F295").Value = 6
$ws.Range("G295").Value = 622.14

$ws.Range("F303").Value = 41
$ws.Range("G303").Value = 8646.49

$ws.Range("B304").Value = 197438.5

$ws.Range("B322").Value = 58047
$ws.Range("D322").Value = 105.54
$ws.Range("E322").Value = 126.1
$ws.Range("F322").Value = 41
$ws.Range("G322").Value = 4327.14

$ws.Range("B323").Value = 47097
$ws.Range("D323").Value = 112.28
$ws.Range("E323").Value = 134.16
$ws.Range("F323").Value = 15
$ws.Range("G323").Value = 1684.2

$ws.Range("F326").Value = 66
$ws.Range("G326").Value = 1962.84

$ws.Range("F328").Value = 62
$ws.Range("G328").Value = 2307.02

$ws.Range("B330").Value = 32341

$ws.Range("F339").Value = 7
$ws.Range("G339").Value = 331.8

$ws.Range("F342").Value = 143
$ws.Range("G342").Value = 4528.81

$ws.Range("F343").Value = 37
$ws.Range("G343").Value = 2662.89

$ws.Range("F345").Value = 80
$ws.Range("G345").Value = 4912.8

$ws.Range("B346").Value = 28789.91

$ws.Range("F353").Value = 17
$ws.Range("G353").Value = 2332.23

$ws.Range("F354").Value = 21
$ws.Range("G354").Value = 1440.39

$ws.Range("B358").Value = 37391.79

$ws.Range("B364").Value = 53602
$ws.Range("E364").Value = 15.69
$ws.Range("F364").Value = -231
$ws.Range("G364").Value = -3037.65

$ws.Range("B365").Value = 65068
$ws.Range("E365").Value = 13.97
$ws.Range("F365").Value = 63
$ws.Range("G365").Value = 828.45

$ws.Range("B366").Value = 65066
$ws.Range("E366").Value = 13.61
$ws.Range("F366").Value = 90
$ws.Range("G366").Value = 1152.9

$ws.Range("B367").Value = 53263
$ws.Range("E367").Value = 15.29
$ws.Range("F367").Value = -309
$ws.Range("G367").Value = -3958.29

$ws.Range("B380").Value = 64925
$ws.Range("E380").Value = 13.97
$ws.Range("F380").Value = 111
$ws.Range("G380").Value = 1459.65

$ws.Range("B381").Value = 45709
$ws.Range("E381").Value = 15.69
$ws.Range("F381").Value = -300
$ws.Range("G381").Value = -3945

$ws.Range("B382").Value = 45702
$ws.Range("E382").Value = 31.43
$ws.Range("F382").Value = -215
$ws.Range("G382").Value = -5654.5

$ws.Range("B383").Value = 64919
$ws.Range("E383").Value = 27.97
$ws.Range("F383").Value = 61
$ws.Range("G383").Value = 1604.3

$ws.Range("F434").Value = 40
$ws.Range("G434").Value = 1305.6

$ws.Range("B435").Value = 1593.98

$ws.Range("B442").Value = 64810
$ws.Range("E442").Value = 291.22
$ws.Range("F442").Value = 5
$ws.Range("G442").Value = 1369.6

$ws.Range("B443").Value = 53319
$ws.Range("E443").Value = 310.64
$ws.Range("F443").Value = -6
$ws.Range("G443").Value = -1643.52

$ws.Range("F453").Value = 27
$ws.Range("G453").Value = 715.77

$ws.Range("B460").Value = 15396.03

$ws.Range("B473").Value = 60022
$ws.Range("E473").Value = 37.22
$ws.Range("F473").Value = -113
$ws.Range("G473").Value = -3709.79

$ws.Range("B474").Value = 64830
$ws.Range("E474").Value = 34.9
$ws.Range("F474").Value = 109
$ws.Range("G474").Value = 3578.47

$ws.Range("F508").Value = 62
$ws.Range("G508").Value = 6444.28

$ws.Range("B510").Value = 26458.9

$ws.Range("F551").Value = 8
$ws.Range("G551").Value = 1145.04

$ws.Range("F552").Value = 28
$ws.Range("G552").Value = 2850.12

$ws.Range("B560").Value = 9289.459999999999

$ws.Range("F575").Value = 5
$ws.Range("G575").Value = 165.3

$ws.Range("B583").Value = 28970.28

$ws.Range("F599").Value = 2070
$ws.Range("G599").Value = 337637.7

$ws.Range("F602").Value = 353
$ws.Range("G602").Value = 51061.45

$ws.Range("B606").Value = 522213.23

$ws.Range("F612").Value = 35
$ws.Range("G612").Value = 1434.65

$ws.Range("B618").Value = 46988.06

$ws.Range("B619").Value = 1996226.33

$ws.Range("B620").Value = 1996226.33

